# TS 5.3 Jatai Sanskrit Corrections
#
# Target: a run containing a lone "x" (red-on-highlight correction mark,
# originally green-highlighted) immediately followed by a run containing
# "iÉiÉÑþÈ |". The "i" that starts the second run actually belongs with
# the correction mark, so it is folded into the first run (x -> xi) and
# dropped from the front of the second run (iÉiÉÑþÈ | -> ÉiÉÑþÈ |). The
# corrected run also switches from a black/green-highlight flag to a
# red/yellow-highlight flag.
#
# The phrase "xiÉiÉÑþÈ" is unique in the whole document, so we can find
# it reliably and then work with precise character offsets relative to
# that match.

$d = $word.ActiveDocument

$anchor = $d.Content
$anchor.Find.Execute("xiÉiÉÑþÈ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$matchStart = $anchor.Start

# Offsets (relative to $matchStart) of the "xi" that must become one
# red/yellow-highlighted run, and of the paragraph-run boundary right
# after "...È |" where an existing (unrelated) run already begins.
$xiStart = $matchStart
$xiEnd = $matchStart + 2
$afterRunEnd = $matchStart + 10

# The run immediately following "...È |" already carries the exact same
# character formatting, so Word's find/replace engine would otherwise
# silently merge it into the run we are about to rewrite. Nudge its
# color off momentarily so the two stay distinct, then restore it.
$guard = $d.Range($afterRunEnd, $afterRunEnd + 1)
$guard.Font.Color = 123456

$target = $d.Range($xiStart, $xiEnd)
$find = $target.Find
$find.ClearFormatting()
$find.Text = "xi"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "xi"
$find.Replacement.Highlight = $true
$find.Replacement.Font.Color = 255
$find.Execute("xi", $true, $false, $false, $false, $false, $true, 1, $false, "xi", 1, $true) | Out-Null

$guard2 = $d.Range($afterRunEnd, $afterRunEnd + 1)
$guard2.Font.Color = 0
